# Generate Report for Handback
#
# Updates the localization-status report after a handback run:
#  - "In Translation" status becomes "Handed back: in sync with en-US" everywhere
#  - the zh-cn / de-de sheets get their "Latest Target File" (I) and
#    "Latest Handback File" (J) columns populated, plus an actual
#    "Latest Handback DateTime" (K) timestamp instead of the zero date
#  - a couple of columns are widened to fit the newly-populated values

$wb = $excel.ActiveWorkbook

$urlForDoc = @{
    "17968538-7582-4e48-8a4f-6ea36eb74e3e.md" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10b080869b6aba882575eb783efeef4b0955dba7/e2e/17968538-7582-4e48-8a4f-6ea36eb74e3e.md";
    "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10b080869b6aba882575eb783efeef4b0955dba7/e2e/8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md";
}

# ---------------------------------------------------------------------------
# 1) Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3)
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback
#    DateTime for both data rows, plus hyperlinks on the new "Target File"
#    cells (rebuild the hyperlink collection so ordering/relationship ids
#    come out the same way Excel lays them out: row-major, left to right).
# ---------------------------------------------------------------------------

$wsZhCn.Range("A2").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlForDoc["17968538-7582-4e48-8a4f-6ea36eb74e3e.md"], $null, $null, "17968538-7582-4e48-8a4f-6ea36eb74e3e.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlForDoc["17968538-7582-4e48-8a4f-6ea36eb74e3e.md"], $null, $null, "17968538-7582-4e48-8a4f-6ea36eb74e3e.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlForDoc["8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md"], $null, $null, "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlForDoc["8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md"], $null, $null, "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md")

$wsZhCn.Range("J2").Value = "17968538-7582-4e48-8a4f-6ea36eb74e3e.af2b41395b0065a2633bd7a6ef51037d5b16d4fa.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.143d23643586dfa1c6bf6ad9ce12627ec431870b.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-25 18:24:58"
$wsZhCn.Range("K3").Value = "2016-08-25 18:24:58"

# ---------------------------------------------------------------------------
# 3) de-de sheet: same shape of update, different handback timestamp.
# ---------------------------------------------------------------------------

$wsDeDe.Range("A2").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlForDoc["17968538-7582-4e48-8a4f-6ea36eb74e3e.md"], $null, $null, "17968538-7582-4e48-8a4f-6ea36eb74e3e.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlForDoc["17968538-7582-4e48-8a4f-6ea36eb74e3e.md"], $null, $null, "17968538-7582-4e48-8a4f-6ea36eb74e3e.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlForDoc["8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md"], $null, $null, "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlForDoc["8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md"], $null, $null, "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.md")

$wsDeDe.Range("J2").Value = "17968538-7582-4e48-8a4f-6ea36eb74e3e.af2b41395b0065a2633bd7a6ef51037d5b16d4fa.de-de.xlf"
$wsDeDe.Range("J3").Value = "8c47fb73-0f65-4ed8-ab4b-23f0e11e75e3.143d23643586dfa1c6bf6ad9ce12627ec431870b.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-25 18:25:15"
$wsDeDe.Range("K3").Value = "2016-08-25 18:25:15"

# ---------------------------------------------------------------------------
# 4) Widen columns that now hold the longer handback file names / dates.
# ---------------------------------------------------------------------------

$wsOverview.Range("E:E").ColumnWidth = 29.166666666666668
$wsOverview.Range("F:F").ColumnWidth = 29.166666666666668

$wsZhCn.Range("C:C").ColumnWidth = 29.166666666666668
$wsZhCn.Range("I:I").ColumnWidth = 39.166666666666664
$wsZhCn.Range("J:J").ColumnWidth = 39.166666666666664

$wsDeDe.Range("C:C").ColumnWidth = 29.166666666666668
$wsDeDe.Range("I:I").ColumnWidth = 39.166666666666664
$wsDeDe.Range("J:J").ColumnWidth = 39.166666666666664
